$wb = $excel.ActiveWorkbook

# --- Rename sheet "longrange_FINAL" -> "longrange" ---
$longrange = $wb.Worksheets.Item("longrange_FINAL")
$longrange.Name = "longrange"

$raw = $wb.Worksheets.Item("RawInjectionData")

# --- New data entered at the bottom of the injection-site table ---
$raw.Range("G146").Value = 643749624
$raw.Range("H146").Value = "S1 L5"

# G147 is emphasised (bold, black Arial) to flag it out from the rest
$raw.Range("A146").Copy()
$raw.Range("G147").PasteSpecial(-4122)  # xlPasteFormats
$raw.Range("G147").Font.Bold = $true
$raw.Range("G147").Value = 266486371
$raw.Range("H147").Value = "S1 L234"

$raw.Range("G148").Value = 100142580
$raw.Range("H148").Value = "CP at S1 terminals (opp hemisphere)"

# --- Update selection/view state ---
# Selection left behind on the "longrange" sheet before switching away from it
$longrange.Range("H3").Select()

# RawInjectionData becomes the active/selected sheet, scrolled near the new rows
$raw.Activate()
$raw.Range("G84").Select()
